# Updated progress log: append a new time-tracking entry (row 59) before
# the "Total" row of the last week block ("Week 11" at A53, rows 54-58),
# pushing the existing Total row down to row 60 and extending its SUM
# range to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row just above the current "Total" row (row 59),
# shifting the Total row (and everything below it) down by one.
$ws.Rows.Item(59).Insert()

# Fill in the new day's data in the freshly inserted row 59.
$ws.Range("A59").Value = 42290                  # Tuesday, October 13, 2015
$ws.Range("B59").Value = 0.375                  # 9:00 AM
$ws.Range("C59").Value = 0.45833333333333331    # 11:00 AM
$ws.Range("D59").Value = 0                      # Break

# Duration formula, matching the pattern used by the rows above it.
$ws.Range("E59").Formula = "=MOD(C59-B59,1)*24-D59"

# Re-assert E54's formula so the shared formula group spans E54:E59.
$ws.Range("E54").Formula = "=MOD(C54-B54,1)*24-D54"

# The old Total row (now row 60) needs its SUM widened to include E59.
$ws.Range("E60").Formula = "=SUM(E54:E59)"

# Update the active selection to mirror the author's cursor position.
$ws.Range("F62").Select() | Out-Null
